$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title block text updates ---
# Police Commissioner name change
$ws.Range("M6").Value = "Jessica S. Tisch"

# Volume/Number text (week 47 -> 48)
$ws.Range("A8").Value = "Volume 31   Number  48"

# Report covering week dates
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Crime statistics table updates (rows 14-30) ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 26
$ws.Range("K15").Value = 50
$ws.Range("N15").Value = -43.478260869565
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -22.222222222222
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -14.814814814814
$ws.Range("I16").Value = 407
$ws.Range("J16").Value = 346
$ws.Range("K16").Value = 17.630057803468
$ws.Range("L16").Value = 11.813186813186
$ws.Range("M16").Value = -7.077625570776
$ws.Range("N16").Value = -74.109414758269
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 58
$ws.Range("G17").Value = 49
$ws.Range("H17").Value = 18.367346938775
$ws.Range("I17").Value = 725
$ws.Range("J17").Value = 614
$ws.Range("K17").Value = 18.078175895765
$ws.Range("L17").Value = 34.758364312267
$ws.Range("M17").Value = 126.5625
$ws.Range("N17").Value = -7.878017789072
$ws.Range("C18").Value = 5
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 209
$ws.Range("K18").Value = 39.333333333333
$ws.Range("L18").Value = 7.731958762886
$ws.Range("M18").Value = -9.523809523809
$ws.Range("N18").Value = -80.41237113402
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 5.263157894736
$ws.Range("I19").Value = 554
$ws.Range("J19").Value = 544
$ws.Range("K19").Value = 1.838235294117
$ws.Range("L19").Value = -0.359712230215
$ws.Range("M19").Value = 41.687979539641
$ws.Range("N19").Value = -49.498632634457
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 6.666666666666
$ws.Range("I20").Value = 183
$ws.Range("J20").Value = 214
$ws.Range("K20").Value = -14.485981308411
$ws.Range("L20").Value = -8.040201005025
$ws.Range("M20").Value = 16.56050955414
$ws.Range("N20").Value = -87.508532423208
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 12.903225806451
$ws.Range("F21").Value = 151
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = 4.861111111111
$ws.Range("I21").Value = 2124
$ws.Range("J21").Value = 1896
$ws.Range("K21").Value = 12.025316455696
$ws.Range("L21").Value = 12.143611404435
$ws.Range("M21").Value = 35.459183673469
$ws.Range("N21").Value = -65.094494658997
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 36.363636363636
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = 25
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = -100
$ws.Range("M23").Value = 8.108108108108
$ws.Range("F24").Value = 169
$ws.Range("H24").Value = 35.2
$ws.Range("I24").Value = 1940
$ws.Range("J24").Value = 1460
$ws.Range("K24").Value = 32.876712328767
$ws.Range("L24").Value = 42.124542124542
$ws.Range("M24").Value = 87.43961352657
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F25").Value = 117
$ws.Range("G25").Value = 57
$ws.Range("H25").Value = 105.263157894737
$ws.Range("I25").Value = 1216
$ws.Range("J25").Value = 617
$ws.Range("K25").Value = 97.08265802269
$ws.Range("L25").Value = 93.939393939393
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = 73.333333333333
$ws.Range("F26").Value = 87
$ws.Range("G26").Value = 74
$ws.Range("H26").Value = 17.567567567567
$ws.Range("I26").Value = 1006
$ws.Range("J26").Value = 907
$ws.Range("K26").Value = 10.915104740904
$ws.Range("L26").Value = 39.722222222222
$ws.Range("M26").Value = 38.950276243093
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 57
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = 54.054054054054
$ws.Range("L27").Value = 18.75
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 112
$ws.Range("J28").Value = 84
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = 45.454545454545
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("F29").Value = 4
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("I29").Value = 25
$ws.Range("K29").Value = 127.272727272727
$ws.Range("L29").Value = -16.666666666666
$ws.Range("M29").Value = -32.432432432432
$ws.Range("N29").Value = -83.333333333333
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("F30").Value = 3
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "***.*"
$ws.Range("I30").Value = 20
$ws.Range("K30").Value = 81.818181818181
$ws.Range("L30").Value = -20
$ws.Range("M30").Value = -20
$ws.Range("N30").Value = -85.507246376811
